$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Description values for StopId..OfficeCode rows (column F) ---
# (written in the same order the shared-string table records them)
$ws.Range("F5").Value  = "Primary id for auto increatment."
$ws.Range("F7").Value  = "For active or inactive record."
$ws.Range("F8").Value  = "Get Login UserId ."
$ws.Range("F9").Value  = "Get Current Datetime when user Insert Record."
$ws.Range("F10").Value = "Keep IP Address of User System."
$ws.Range("F11").Value = "Get OfficeCode  (MstOIS) in Numaric "
$ws.Range("F6").Value  = "Take StopName as String."

# --- Apply the new "description" font (Times New Roman 11, dark grey FF1F1F1F) ---
# and left alignment to the newly-filled description cells
$descRange = $ws.Range("F5:F11")
$descRange.Font.Name = "Times New Roman"
$descRange.Font.Family = 1
$descRange.Font.Color = 2039583
$descRange.HorizontalAlignment = -4131

# --- Highlight the OfficeCode row's column-name cell with a green fill ---
$ws.Range("B11").Interior.Color = 5296274

# --- Widen column F so the long description text is readable ---
$ws.Columns("F").ColumnWidth = 50.9

# --- Move the active selection (matches the saved cursor position) ---
$ws.Range("D12").Select()
